$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.991.54"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.641.68"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'215.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").Value = "'19.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").Value = "'0.0796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "1.869.81"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "1.657.86"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").Value = "'63.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").Value = "26.018.13"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "'194.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "'4.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").Value = "'9.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -1.10%  "
$ws.Range("E24").Value = "  +4.68%  "
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "'142.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("E35").Value = "  +1.47%  "
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").Value = "1.130.99"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("E38").Value = "  -0.85%  "
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").Value = "'5.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("D42").Value = "'99.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").Value = "'0.797"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  +2.85%  "
$ws.Range("D45").Value = "'56.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("E46").Value = "  +3.54%  "
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("D48").Value = "'7.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'0.0951"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.65%  "